$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new row's cells as Text first so that date-looking strings
# (e.g. "2025-06-19") are stored as literal text instead of being
# auto-converted into a date serial number by Excel's input parsing.
$ws.Range("A19:F19").NumberFormat = "@"

$ws.Range("A19").Value = "edit1"
$ws.Range("B19").Value = "riya-morankar"
$ws.Range("C19").Value = "Merged"
$ws.Range("D19").Value = "N/A"
$ws.Range("E19").Value = "2025-06-19"
$ws.Range("F19").Value = "8b59deccfda5d1e814ac4d6141e02dc3d5f12b19"

# Restore the default "Normal" style so the new cells match the rest of
# the sheet (which uses the workbook's default style with no explicit
# per-cell style index), while keeping the values stored as text.
$ws.Range("A19:F19").Style = "Normal"
